$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing column B ("Jun_13") so that:
#   old column B (Jun_13) shifts to column D
#   old column C (Jun_10) shifts to column E
$ws.Range("B1:C1").EntireColumn.Insert()

# New column headers (dates)
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"

# Give the two newly inserted columns the same width as the neighboring
# rating columns (matches the original ~8.0 character width).
$ws.Range("C1:D1").EntireColumn.ColumnWidth = 7.1666666666666667

# Fill the new rating columns with the default "UN" (unchanged) rating for
# every analyst row, mirroring the existing Jun_13/Jun_10 columns.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
